# Update NATMI LR-pair sheet (Ccl28-Ackr2) with new TPM-based results.
# A new "Sending cluster" (ECs) row is introduced; it sorts before the
# existing FAPs / MuSCs rows, so the sheet now has 3 data rows instead of 2,
# and the numeric columns for all three rows reflect the refreshed TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Sending cluster = ECs --------------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05794866666666667
$ws.Range("H2").Value = 0.173846
$ws.Range("I2").Value = 0.235800444619869
$ws.Range("J2").Value = 0.235800444619869
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.51723133333333
$ws.Range("N2").Value = 34.551694
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.6674081994582223
$ws.Range("R2").Value = 6.006673795124
$ws.Range("S2").Value = 0.235800444619869
$ws.Range("T2").Value = 0.235800444619869

# ---- Row 3: Sending cluster = FAPs -------------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09370299999999999
$ws.Range("H3").Value = 0.281109
$ws.Range("I3").Value = 0.3812893433650861
$ws.Range("J3").Value = 0.3812893433650861
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.51723133333333
$ws.Range("N3").Value = 34.551694
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.079199127627333
$ws.Range("R3").Value = 9.712792148645999
$ws.Range("S3").Value = 0.3812893433650861
$ws.Range("T3").Value = 0.3812893433650861

# ---- Row 4 (new): Sending cluster = MuSCs ------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ccl28"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.09410133333333333
$ws.Range("H4").Value = 0.282304
$ws.Range("I4").Value = 0.3829102120150449
$ws.Range("J4").Value = 0.3829102120150449
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.51723133333333
$ws.Range("N4").Value = 34.551694
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.083786824775111
$ws.Range("R4").Value = 9.754081422975998
$ws.Range("S4").Value = 0.3829102120150449
$ws.Range("T4").Value = 0.3829102120150449
